# Fix contact information missing from short resumes.
#
# The "short" resume template was missing the contact-info line that
# appears right under the candidate's name. Re-insert it as its own
# centered paragraph immediately after the "Dheeraj Chand" heading
# paragraph and before "PROFESSIONAL SUMMARY".

$d = $word.ActiveDocument

$contactInfo = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

# Use Find/Replace with a literal paragraph-mark ("^p") in the replacement
# text so Word splits the "Dheeraj Chand" paragraph into two paragraphs:
# the original (untouched formatting) plus a brand-new plain paragraph
# that inherits the same centered alignment but carries no extra run
# formatting (matches the target diff exactly - no stray <w:rPr>).
$d.Content.Find.Execute("Dheeraj Chand", $false, $false, $false, $false, `
                         $false, $true, 1, $false, `
                         "Dheeraj Chand^p$contactInfo", 2)
